# Horarios actualizados Linea 141 - 638
# Refresh scrape: new scrape time 09:35:26, updated row counts, re-sorted by Hora_Llegada
$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2,1).Value = "Última actualización: 09:35:26"
$ws1.Cells.Item(3,1).Value = "Total filas: 137"
$ws1Data = @(
    @(6, "03:45:25", "03:47", "14_ABASTO", 2, "LP1912"),
    @(7, "03:45:25", "04:01", "81_EL PELIGRO", 16, "LP1912"),
    @(8, "03:45:25", "04:46", "215A_EL PATO", 61, "LP1912"),
    @(9, "03:45:25", "04:53", "11_ETCHEVERRY", 68, "LP1912"),
    @(10, "04:56:49", "05:13", "14_ABASTO", 17, "LP1912"),
    @(11, "03:45:25", "05:16", "17_ROMERO", 91, "LP1912"),
    @(12, "04:45:05", "05:16", "14_ABASTO", 31, "LP1912"),
    @(13, "03:45:25", "05:22", "23_HERNANDEZ", 97, "LP1912"),
    @(14, "05:26:08", "05:28", "14_ABASTO", 2, "LP1912"),
    @(15, "03:45:25", "05:34", "215B_EL PATO", 109, "LP1912"),
    @(16, "04:18:02", "05:34", "14_ABASTO", 76, "LP1912"),
    @(17, "04:18:02", "05:35", "215B_EL PATO", 77, "LP1912"),
    @(18, "03:45:25", "05:37", "14_ABASTO", 112, "LP1912"),
    @(19, "04:18:02", "05:46", "15_ABASTO", 88, "LP1912"),
    @(20, "04:45:05", "06:04", "16_SANTA ANA", 79, "LP1912"),
    @(21, "04:18:02", "06:05", "16_SANTA ANA", 107, "LP1912"),
    @(22, "04:56:49", "06:11", "215A_EL PATO", 75, "LP1912"),
    @(23, "04:18:02", "06:12", "215A_EL PATO", 114, "LP1912"),
    @(24, "04:18:02", "06:14", "225_HARAS DEL SUR", 116, "LP1912"),
    @(25, "04:45:05", "06:21", "26_HERNANDEZ", 96, "LP1912"),
    @(26, "06:25:43", "06:26", "86_EST CHICA-ESC AGRARIA", 1, "LP1912"),
    @(27, "04:45:05", "06:27", "23_HERNANDEZ", 102, "LP1912"),
    @(28, "06:25:43", "06:28", "23_HERNANDEZ", 3, "LP1912"),
    @(29, "04:56:49", "06:29", "86_EST CHICA-ESC AGRARIA", 93, "LP1912"),
    @(30, "04:45:05", "06:30", "86_EST CHICA-ESC AGRARIA", 105, "LP1912"),
    @(31, "04:45:05", "06:31", "16_SANTA ANA", 106, "LP1912"),
    @(32, "04:45:05", "06:44", "225_C ROCA-H SUR", 119, "LP1912"),
    @(33, "05:55:25", "06:44", "26_HERNANDEZ", 49, "LP1912"),
    @(34, "04:56:49", "06:46", "215C_EL PATO", 110, "LP1912"),
    @(35, "05:26:08", "06:47", "215C_EL PATO", 81, "LP1912"),
    @(36, "05:55:25", "06:59", "14_ABASTO", 64, "LP1912"),
    @(37, "05:26:08", "07:00", "14_ABASTO", 94, "LP1912"),
    @(38, "06:25:43", "07:01", "16_SANTA ANA", 36, "LP1912"),
    @(39, "05:55:25", "07:04", "23_HERNANDEZ", 69, "LP1912"),
    @(40, "05:26:08", "07:05", "15_ABASTO", 99, "LP1912"),
    @(41, "05:26:08", "07:05", "23_HERNANDEZ", 99, "LP1912"),
    @(42, "05:26:08", "07:06", "10_OLMOS", 100, "LP1912"),
    @(43, "05:26:08", "07:07", "225_GOMEZ", 101, "LP1912"),
    @(44, "05:26:08", "07:11", "215A_EL PATO", 105, "LP1912"),
    @(45, "06:55:02", "07:12", "215A_EL PATO", 17, "LP1912"),
    @(46, "06:25:43", "07:14", "26_HERNANDEZ", 49, "LP1912"),
    @(47, "05:55:25", "07:15", "11_ETCHEVERRY", 80, "LP1912"),
    @(48, "05:26:08", "07:16", "11_ETCHEVERRY", 110, "LP1912"),
    @(49, "06:55:02", "07:17", "16_SANTA ANA", 22, "LP1912"),
    @(50, "05:26:08", "07:21", "26_HERNANDEZ", 115, "LP1912"),
    @(51, "05:26:08", "07:23", "10_OLMOS", 117, "LP1912"),
    @(52, "05:55:25", "07:30", "10_OLMOS", 95, "LP1912"),
    @(53, "05:55:25", "07:31", "16_SANTA ANA", 96, "LP1912"),
    @(54, "05:55:25", "07:31", "11_ETCHEVERRY", 96, "LP1912"),
    @(55, "06:55:02", "07:32", "11_ETCHEVERRY", 37, "LP1912"),
    @(56, "06:55:02", "07:32", "16_SANTA ANA", 37, "LP1912"),
    @(57, "05:55:25", "07:32", "84_COLONIA URQUIZA-ESC 49", 97, "LP1912"),
    @(58, "07:19:29", "07:35", "23_HERNANDEZ", 16, "LP1912"),
    @(59, "05:55:25", "07:36", "27_EL RETIRO", 101, "LP1912"),
    @(60, "06:55:02", "07:37", "27_EL RETIRO", 42, "LP1912"),
    @(61, "05:55:25", "07:39", "10_OLMOS", 104, "LP1912"),
    @(62, "07:19:29", "07:46", "16_SANTA ANA", 27, "LP1912"),
    @(63, "05:55:25", "07:47", "14_ABASTO", 112, "LP1912"),
    @(64, "06:55:02", "07:48", "14_ABASTO", 53, "LP1912"),
    @(65, "07:50:16", "07:50", "10_OLMOS", 0, "LP1912"),
    @(66, "05:55:25", "07:51", "215D_EL PATO", 116, "LP1912"),
    @(67, "06:55:02", "07:52", "215D_EL PATO", 57, "LP1912"),
    @(68, "07:19:29", "07:59", "23_HERNANDEZ", 40, "LP1912"),
    @(69, "06:25:43", "08:01", "23_HERNANDEZ", 96, "LP1912"),
    @(70, "07:19:29", "08:03", "11_ETCHEVERRY", 44, "LP1912"),
    @(71, "06:55:02", "08:03", "23_HERNANDEZ", 68, "LP1912"),
    @(72, "08:02:22", "08:05", "23_HERNANDEZ", 3, "LP1912"),
    @(73, "07:19:29", "08:10", "16_SANTA ANA", 51, "LP1912"),
    @(74, "07:50:16", "08:11", "16_SANTA ANA", 21, "LP1912"),
    @(75, "06:25:43", "08:12", "15_ABASTO", 107, "LP1912"),
    @(76, "07:50:16", "08:13", "10_OLMOS", 23, "LP1912"),
    @(77, "06:55:02", "08:21", "26_HERNANDEZ", 86, "LP1912"),
    @(78, "06:25:43", "08:22", "16_P MOR-SANTA ANA", 117, "LP1912"),
    @(79, "06:55:02", "08:23", "16_P MOR-SANTA ANA", 88, "LP1912"),
    @(80, "06:25:43", "08:23", "215B_EL PATO", 118, "LP1912"),
    @(81, "06:55:02", "08:27", "84_COLONIA URQUIZA-ESC 49", 92, "LP1912"),
    @(82, "07:50:16", "08:30", "23_HERNANDEZ", 40, "LP1912"),
    @(83, "08:02:22", "08:33", "10_OLMOS", 31, "LP1912"),
    @(84, "08:02:22", "08:34", "23_HERNANDEZ", 32, "LP1912"),
    @(85, "08:32:09", "08:37", "23_HERNANDEZ", 5, "LP1912"),
    @(86, "06:55:02", "08:42", "81_EL PELIGRO", 107, "LP1912"),
    @(87, "07:19:29", "08:43", "14_ABASTO", 84, "LP1912"),
    @(88, "07:50:16", "08:44", "14_ABASTO", 54, "LP1912"),
    @(89, "08:32:09", "08:53", "10_OLMOS", 21, "LP1912"),
    @(90, "06:55:02", "08:54", "17_ROMERO", 119, "LP1912"),
    @(91, "07:19:29", "09:01", "215A_EL PATO", 102, "LP1912"),
    @(92, "07:50:16", "09:02", "215A_EL PATO", 72, "LP1912"),
    @(93, "08:02:22", "09:03", "11_ETCHEVERRY", 61, "LP1912"),
    @(94, "08:32:09", "09:04", "11_ETCHEVERRY", 32, "LP1912"),
    @(95, "08:32:09", "09:05", "23_HERNANDEZ", 33, "LP1912"),
    @(96, "07:19:29", "09:10", "16_P MOR-SANTA ANA", 111, "LP1912"),
    @(97, "07:50:16", "09:11", "16_P MOR-SANTA ANA", 81, "LP1912"),
    @(98, "08:32:09", "09:13", "10_OLMOS", 41, "LP1912"),
    @(99, "07:19:29", "09:16", "27_EL RETIRO", 117, "LP1912"),
    @(100, "07:50:16", "09:17", "27_EL RETIRO", 87, "LP1912"),
    @(101, "07:50:16", "09:21", "26_HERNANDEZ", 91, "LP1912"),
    @(102, "08:02:22", "09:22", "16_SANTA ANA", 80, "LP1912"),
    @(103, "08:02:22", "09:23", "11_ETCHEVERRY", 81, "LP1912"),
    @(104, "08:32:09", "09:23", "16_SANTA ANA", 51, "LP1912"),
    @(105, "07:50:16", "09:23", "17_ROMERO", 93, "LP1912"),
    @(106, "07:50:16", "09:24", "11_ETCHEVERRY", 94, "LP1912"),
    @(107, "07:50:16", "09:28", "16_SANTA ANA", 98, "LP1912"),
    @(108, "07:50:16", "09:32", "15_ABASTO", 102, "LP1912"),
    @(109, "07:50:16", "09:33", "10_OLMOS", 103, "LP1912"),
    @(110, "08:56:29", "09:34", "23_HERNANDEZ", 38, "LP1912"),
    @(111, "08:56:29", "09:34", "16_SANTA ANA", 38, "LP1912"),
    @(112, "08:32:09", "09:35", "16_SANTA ANA", 63, "LP1912"),
    @(113, "08:48:08", "09:35", "23_HERNANDEZ", 47, "LP1912"),
    @(114, "09:35:26", "09:39", "23_HERNANDEZ", 4, "LP1912"),
    @(115, "07:50:16", "09:42", "215C_EL PATO", 112, "LP1912"),
    @(116, "08:02:22", "09:43", "14_ABASTO", 101, "LP1912"),
    @(117, "07:50:16", "09:44", "14_ABASTO", 114, "LP1912"),
    @(118, "09:35:26", "09:46", "16_SANTA ANA", 11, "LP1912"),
    @(119, "08:32:09", "09:52", "15_ABASTO", 80, "LP1912"),
    @(120, "08:56:29", "09:53", "10_OLMOS", 57, "LP1912"),
    @(121, "09:35:26", "09:58", "16_SANTA ANA", 23, "LP1912"),
    @(122, "09:35:26", "10:03", "11_ETCHEVERRY", 28, "LP1912"),
    @(123, "08:56:29", "10:10", "16_P MOR-SANTA ANA", 74, "LP1912"),
    @(124, "08:32:09", "10:11", "16_P MOR-SANTA ANA", 99, "LP1912"),
    @(125, "09:35:26", "10:12", "15_ABASTO", 37, "LP1912"),
    @(126, "09:35:26", "10:13", "10_OLMOS", 38, "LP1912"),
    @(127, "08:32:09", "10:21", "26_HERNANDEZ", 109, "LP1912"),
    @(128, "08:32:09", "10:22", "17_ROMERO", 110, "LP1912"),
    @(129, "09:35:26", "10:23", "11_ETCHEVERRY", 48, "LP1912"),
    @(130, "08:56:29", "10:26", "215A_EL PATO", 90, "LP1912"),
    @(131, "08:32:09", "10:27", "215A_EL PATO", 115, "LP1912"),
    @(132, "09:35:26", "10:34", "23_HERNANDEZ", 59, "LP1912"),
    @(133, "08:48:08", "10:42", "17_ROMERO", 114, "LP1912"),
    @(134, "08:56:29", "10:43", "14_ABASTO", 107, "LP1912"),
    @(135, "08:48:08", "10:44", "14_ABASTO", 116, "LP1912"),
    @(136, "09:35:26", "10:54", "27_EL RETIRO", 79, "LP1912"),
    @(137, "09:35:26", "11:02", "215C_EL PATO", 87, "LP1912"),
    @(138, "09:35:26", "11:06", "16_P MOR-167 Y 521", 91, "LP1912"),
    @(139, "09:35:26", "11:19", "86_EST CHICA-ESC AGRARIA", 104, "LP1912"),
    @(140, "09:35:26", "11:21", "26_HERNANDEZ", 106, "LP1912"),
    @(141, "09:35:26", "11:27", "225_C ROCA-H SUR", 112, "LP1912"),
    @(142, "09:35:26", "11:32", "81_EL PELIGRO", 117, "LP1912")
)
foreach ($row in $ws1Data) {
    $r = $row[0]
    $ws1.Cells.Item($r,1).Value = $row[1]
    $ws1.Cells.Item($r,2).Value = $row[2]
    $ws1.Cells.Item($r,3).Value = $row[3]
    $ws1.Cells.Item($r,4).Value = $row[4]
    $ws1.Cells.Item($r,5).Value = $row[5]
}

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2,1).Value = "Última actualización: 09:35:26"
$ws2.Cells.Item(3,1).Value = "Total filas: 18"
$ws2Data = @(
    @(6, "03:45:25", "04:46", "215A_EL PATO", 61, "LP1912"),
    @(7, "03:45:25", "05:34", "215B_EL PATO", 109, "LP1912"),
    @(8, "04:18:02", "05:35", "215B_EL PATO", 77, "LP1912"),
    @(9, "04:56:49", "06:11", "215A_EL PATO", 75, "LP1912"),
    @(10, "04:18:02", "06:12", "215A_EL PATO", 114, "LP1912"),
    @(11, "04:56:49", "06:46", "215C_EL PATO", 110, "LP1912"),
    @(12, "05:26:08", "06:47", "215C_EL PATO", 81, "LP1912"),
    @(13, "05:26:08", "07:11", "215A_EL PATO", 105, "LP1912"),
    @(14, "06:55:02", "07:12", "215A_EL PATO", 17, "LP1912"),
    @(15, "05:55:25", "07:51", "215D_EL PATO", 116, "LP1912"),
    @(16, "06:55:02", "07:52", "215D_EL PATO", 57, "LP1912"),
    @(17, "06:25:43", "08:23", "215B_EL PATO", 118, "LP1912"),
    @(18, "07:19:29", "09:01", "215A_EL PATO", 102, "LP1912"),
    @(19, "07:50:16", "09:02", "215A_EL PATO", 72, "LP1912"),
    @(20, "07:50:16", "09:42", "215C_EL PATO", 112, "LP1912"),
    @(21, "08:56:29", "10:26", "215A_EL PATO", 90, "LP1912"),
    @(22, "08:32:09", "10:27", "215A_EL PATO", 115, "LP1912"),
    @(23, "09:35:26", "11:02", "215C_EL PATO", 87, "LP1912")
)
foreach ($row in $ws2Data) {
    $r = $row[0]
    $ws2.Cells.Item($r,1).Value = $row[1]
    $ws2.Cells.Item($r,2).Value = $row[2]
    $ws2.Cells.Item($r,3).Value = $row[3]
    $ws2.Cells.Item($r,4).Value = $row[4]
    $ws2.Cells.Item($r,5).Value = $row[5]
}

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2,1).Value = "Última actualización: 09:35:26"
$ws3.Cells.Item(3,1).Value = "Total filas: 26"
$ws3Data = @(
    @(6, "04:56:49", "05:43", "215A_LA PLATA", 47, "L6173"),
    @(7, "03:45:25", "05:44", "215A_LA PLATA", 119, "L6173"),
    @(8, "04:56:49", "06:08", "215A_LA PLATA", 72, "L6173"),
    @(9, "04:18:02", "06:09", "215A_LA PLATA", 111, "L6173"),
    @(10, "04:56:49", "06:32", "215C_LA PLATA", 96, "L6203"),
    @(11, "04:45:05", "06:33", "215C_LA PLATA", 108, "L6203"),
    @(12, "06:25:43", "06:59", "215B_LP-P MOR-1 Y 57", 34, "L6173"),
    @(13, "05:26:08", "07:00", "215B_LP-P MOR-1 Y 57", 94, "L6173"),
    @(14, "05:55:25", "07:35", "215A_LA PLATA", 100, "L6173"),
    @(15, "06:25:43", "07:39", "215A_LA PLATA", 74, "L6173"),
    @(16, "06:55:02", "07:42", "215A_LA PLATA", 47, "L6173"),
    @(17, "07:19:29", "07:46", "215A_LA PLATA", 27, "L6173"),
    @(18, "07:50:16", "07:51", "215A_LA PLATA", 1, "L6173"),
    @(19, "06:25:43", "08:06", "215C_LA PLATA", 101, "L6203"),
    @(20, "06:55:02", "08:07", "215C_LA PLATA", 72, "L6203"),
    @(21, "07:19:29", "08:21", "215C_LA PLATA", 62, "L6203"),
    @(22, "08:02:22", "08:25", "215C_LA PLATA", 23, "L6203"),
    @(23, "07:50:16", "08:27", "215C_LA PLATA", 37, "L6203"),
    @(24, "07:19:29", "08:35", "215A_LA PLATA", 76, "L6173"),
    @(25, "06:55:02", "08:36", "215A_LA PLATA", 101, "L6173"),
    @(26, "08:02:22", "08:39", "215A_LA PLATA", 37, "L6173"),
    @(27, "08:32:09", "08:42", "215A_LA PLATA", 10, "L6173"),
    @(28, "07:19:29", "09:09", "215D_LA PLATA", 110, "L6203"),
    @(29, "08:32:09", "10:03", "215B_LP-P MOR-40 Y 115", 91, "L6173"),
    @(30, "08:56:29", "10:54", "215A_LA PLATA", 118, "L6173"),
    @(31, "09:35:26", "11:14", "215C_LA PLATA", 99, "L6203")
)
foreach ($row in $ws3Data) {
    $r = $row[0]
    $ws3.Cells.Item($r,1).Value = $row[1]
    $ws3.Cells.Item($r,2).Value = $row[2]
    $ws3.Cells.Item($r,3).Value = $row[3]
    $ws3.Cells.Item($r,4).Value = $row[4]
    $ws3.Cells.Item($r,5).Value = $row[5]
}
